# Add "continuous frames" annotation columns (G:K) to the N2_HD_list datalist.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Build the new Arial/10pt style once on a scratch cell, then fan it
#        out to the whole G1:K11 block via copy/paste-special (formats only).
#        Doing it this way (rather than touching every cell's Font object)
#        keeps the workbook's cellXfs table minimal -- one new cell format
#        gets created instead of one per cell.
$scratch = $ws.Range("Z1")
$scratch.Font.Name = "Arial"
$scratch.Font.Size = 10
$scratch.Copy()
$ws.Range("G1:K11").PasteSpecial(-4122)   # xlPasteFormats
$scratch.Clear()

# --- 2. Write the continuous-frame ranges, row by row, left to right so the
#        shared-string table is populated in the same order as the source.
$ws.Range("G1").Value = "0-28860"
$ws.Range("H1").Value = "28980-30960"

$ws.Range("G2").Value = "0-8430"
$ws.Range("H2").Value = "8500-end"

$ws.Range("G3").Value = "70-1920"
$ws.Range("H3").Value = "2000-31540"

$ws.Range("G4").Value = "0-7100"
$ws.Range("H4").Value = "7180-13880"
$ws.Range("I4").Value = "13950-31800"

$ws.Range("G5").Value = "120-980"
$ws.Range("H5").Value = "1070-29530"
$ws.Range("I5").Value = "29570-end"

$ws.Range("G6").Value = "100-10100"
$ws.Range("H6").Value = "10270-20690"
$ws.Range("I6").Value = "20800-end"

$ws.Range("G7").Value = "0-280"
$ws.Range("H7").Value = "330-22120"
$ws.Range("I7").Value = "22300-27580"
$ws.Range("J7").Value = "27640-end"

$ws.Range("G8").Value = "0-8020"
$ws.Range("H8").Value = "8120-13000"
$ws.Range("I8").Value = "13040-20730"
$ws.Range("J8").Value = "20950-24780"

$ws.Range("G9").Value = "0-23310"
$ws.Range("H9").Value = "23460-end"

$ws.Range("G10").Value = "0-1330"
$ws.Range("H10").Value = "1370-14260"
$ws.Range("I10").Value = "14380-25520"
$ws.Range("J10").Value = "25680-31030"
$ws.Range("K10").Value = "31100-end"

$ws.Range("G11").Value = "60-19500"
$ws.Range("H11").Value = "19750-26570"
$ws.Range("I11").Value = "26760-29710"
$ws.Range("J11").Value = "29830-end"

# --- 3. Move the active selection to D20, matching the saved view state.
$ws.Range("D20").Select()
